$d = $word.ActiveDocument

# Locate the bullet paragraph that currently reads:
#   "Open API specification (Swagger tool to implement those)"
# This is the last bullet under the "Documentation" heading.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Open API specification*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Open API specification' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Two new bullets need to be inserted right before it:
#   "Open API specification (Swagger tool to implement those)"
#   "Changelogs"
# and the original bullet's text becomes "Proper commenting" (formatting,
# numbering and the trailing run are left untouched).
$target.Range.InsertParagraphBefore()
$target.Range.InsertParagraphBefore()

$p1 = $d.Paragraphs.Item($targetIndex)
$p2 = $d.Paragraphs.Item($targetIndex + 1)
$p3 = $d.Paragraphs.Item($targetIndex + 2)

$p1.Range.Text = "Open API specification (Swagger tool to implement those)"
$p1.Format.SpaceBefore = 0
$p1.Format.SpaceBeforeAuto = 0
$p1.Format.SpaceAfter = 0
$p1.Format.SpaceAfterAuto = 0

$p2.Range.Text = "Changelogs"
$p2.Format.SpaceBefore = 0
$p2.Format.SpaceBeforeAuto = 0
$p2.Format.SpaceAfter = 0
$p2.Format.SpaceAfterAuto = 0

$p3.Range.Text = "Proper commenting"
